$wb = $excel.ActiveWorkbook

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ARM.Range("H32").Value = 17456.756
$ws_ARM.Range("I32").Value = 18241.441
$ws_ARM.Range("J32").Value = 11833.167
$ws_ARM.Range("K32").Value = 18241.441
$ws_ARM.Range("L32").Value = 11833.167
$ws_ARM.Range("M32").Value = -17954.441
$ws_ARM.Range("N32").Value = -12407.167
$ws_ARM.Range("H61").Value = 1565.7222
$ws_ARM.Range("I61").Value = 1615
$ws_ARM.Range("J61").Value = 1393.25
$ws_ARM.Range("K61").Value = 1615
$ws_ARM.Range("L61").Value = 1393.25
$ws_ARM.Range("M61").Value = -1403
$ws_ARM.Range("N61").Value = -1817.25
$ws_ARM.Range("H74").Value = 1310.9546
$ws_ARM.Range("I74").Value = 1094.7
$ws_ARM.Range("J74").Value = 1491.1666
$ws_ARM.Range("K74").Value = 1094.7
$ws_ARM.Range("L74").Value = 1491.1666
$ws_ARM.Range("M74").Value = -220.7
$ws_ARM.Range("N74").Value = -3239.1666
$ws_ARM.Range("H77").Value = 1310.9546
$ws_ARM.Range("I77").Value = 1094.7
$ws_ARM.Range("J77").Value = 1491.1666
$ws_ARM.Range("K77").Value = 5473.5
$ws_ARM.Range("L77").Value = 7455.833000000001
$ws_ARM.Range("M77").Value = -1105.5
$ws_ARM.Range("N77").Value = -16191.833
$ws_ARM.Range("H122").Value = 1748.5652
$ws_ARM.Range("I122").Value = 1801.7894
$ws_ARM.Range("J122").Value = 1495.75
$ws_ARM.Range("K122").Value = 5405.3682
$ws_ARM.Range("L122").Value = 4487.25
$ws_ARM.Range("M122").Value = -2955.3682
$ws_ARM.Range("N122").Value = -9387.25
$ws_ARM.Range("H133").Value = 40646.375
$ws_ARM.Range("J133").Value = 40646.375
$ws_ARM.Range("L133").Value = 40646.375
$ws_ARM.Range("N133").Value = -45706.375
$ws_ARM.Range("H135").Value = 29692.555
$ws_ARM.Range("J135").Value = 29692.555
$ws_ARM.Range("L135").Value = 29692.555
$ws_ARM.Range("N135").Value = -39832.555
$ws_ARM.Range("H136").Value = 1565.7222
$ws_ARM.Range("I136").Value = 1615
$ws_ARM.Range("J136").Value = 1393.25
$ws_ARM.Range("K136").Value = 4845
$ws_ARM.Range("L136").Value = 4179.75
$ws_ARM.Range("M136").Value = -2295
$ws_ARM.Range("N136").Value = -9279.75
$ws_BSM.Range("H94").Value = 781.5599999999999
$ws_BSM.Range("I94").Value = 617.7954999999999
$ws_BSM.Range("K94").Value = 617.7954999999999
$ws_BSM.Range("M94").Value = -166.7954999999999
$ws_BSM.Range("H99").Value = 2115.7896
$ws_BSM.Range("I99").Value = 1166.6666
$ws_BSM.Range("K99").Value = 1166.6666
$ws_BSM.Range("M99").Value = 331.3334
$ws_BSM.Range("H138").Value = 33365.453
$ws_BSM.Range("J138").Value = 33365.453
$ws_BSM.Range("L138").Value = 33365.453
$ws_BSM.Range("N138").Value = -43645.453
$ws_CRP.Range("H16").Value = 1409
$ws_CRP.Range("I16").Value = 1277.5555
$ws_CRP.Range("J16").Value = 2000.5
$ws_CRP.Range("K16").Value = 1277.5555
$ws_CRP.Range("L16").Value = 2000.5
$ws_CRP.Range("M16").Value = -990.5554999999999
$ws_CRP.Range("N16").Value = -2574.5
$ws_CRP.Range("H86").Value = 2743.2727
$ws_CRP.Range("I86").Value = 2710.75
$ws_CRP.Range("J86").Value = 2830
$ws_CRP.Range("K86").Value = 2710.75
$ws_CRP.Range("L86").Value = 2830
$ws_CRP.Range("M86").Value = -1587.75
$ws_CRP.Range("N86").Value = -5076
$ws_CRP.Range("H89").Value = 2743.2727
$ws_CRP.Range("I89").Value = 2710.75
$ws_CRP.Range("J89").Value = 2830
$ws_CRP.Range("K89").Value = 13553.75
$ws_CRP.Range("L89").Value = 14150
$ws_CRP.Range("M89").Value = -7937.75
$ws_CRP.Range("N89").Value = -25382
$ws_CRP.Range("H99").Value = 2206.3809
$ws_CRP.Range("I99").Value = 2464.375
$ws_CRP.Range("J99").Value = 1380.8
$ws_CRP.Range("K99").Value = 2464.375
$ws_CRP.Range("L99").Value = 1380.8
$ws_CRP.Range("M99").Value = -966.375
$ws_CRP.Range("N99").Value = -4376.8
$ws_CRP.Range("H113").Value = 1409
$ws_CRP.Range("I113").Value = 1277.5555
$ws_CRP.Range("J113").Value = 2000.5
$ws_CRP.Range("K113").Value = 1277.5555
$ws_CRP.Range("L113").Value = 2000.5
$ws_CRP.Range("M113").Value = 892.4445000000001
$ws_CRP.Range("N113").Value = -6340.5
$ws_CRP.Range("H126").Value = 2206.3809
$ws_CRP.Range("I126").Value = 2464.375
$ws_CRP.Range("J126").Value = 1380.8
$ws_CRP.Range("K126").Value = 7393.125
$ws_CRP.Range("L126").Value = 4142.4
$ws_CRP.Range("M126").Value = -4923.125
$ws_CRP.Range("N126").Value = -9082.4
$ws_CUL.Range("H5").Value = 1130.2084
$ws_CUL.Range("I5").Value = 1648.625
$ws_CUL.Range("J5").Value = 871
$ws_CUL.Range("K5").Value = 4945.875
$ws_CUL.Range("L5").Value = 2613
$ws_CUL.Range("M5").Value = -4833.875
$ws_CUL.Range("N5").Value = -2837
$ws_CUL.Range("H135").Value = 1130.2084
$ws_CUL.Range("I135").Value = 1648.625
$ws_CUL.Range("J135").Value = 871
$ws_CUL.Range("K135").Value = 14837.625
$ws_CUL.Range("L135").Value = 7839
$ws_CUL.Range("M135").Value = -12302.625
$ws_CUL.Range("N135").Value = -12909
$ws_GSM.Range("H70").Value = 6218
$ws_GSM.Range("J70").Value = 6389.778
$ws_GSM.Range("L70").Value = 6389.778
$ws_GSM.Range("N70").Value = -6929.778
$ws_GSM.Range("H73").Value = 6218
$ws_GSM.Range("J73").Value = 6389.778
$ws_GSM.Range("L73").Value = 6389.778
$ws_GSM.Range("N73").Value = -8261.778
$ws_GSM.Range("H80").Value = 4015
$ws_GSM.Range("I80").Value = 4167.5
$ws_GSM.Range("J80").Value = 3100
$ws_GSM.Range("K80").Value = 4167.5
$ws_GSM.Range("L80").Value = 3100
$ws_GSM.Range("M80").Value = -3169.5
$ws_GSM.Range("N80").Value = -5096
$ws_GSM.Range("H83").Value = 4015
$ws_GSM.Range("I83").Value = 4167.5
$ws_GSM.Range("J83").Value = 3100
$ws_GSM.Range("K83").Value = 20837.5
$ws_GSM.Range("L83").Value = 15500
$ws_GSM.Range("M83").Value = -15845.5
$ws_GSM.Range("N83").Value = -25484
$ws_GSM.Range("H122").Value = 3566.9565
$ws_GSM.Range("I122").Value = 3802.7778
$ws_GSM.Range("K122").Value = 11408.3334
$ws_GSM.Range("M122").Value = -8958.3334
$ws_LTW.Range("H68").Value = 2535.4285
$ws_LTW.Range("I68").Value = 1370.7142
$ws_LTW.Range("K68").Value = 1370.7142
$ws_LTW.Range("M68").Value = -621.7141999999999
$ws_LTW.Range("H71").Value = 2535.4285
$ws_LTW.Range("I71").Value = 1370.7142
$ws_LTW.Range("K71").Value = 6853.571
$ws_LTW.Range("M71").Value = -3109.571
$ws_LTW.Range("H122").Value = 11368591
$ws_LTW.Range("I122").Value = 17862892
$ws_LTW.Range("J122").Value = 3563.125
$ws_LTW.Range("K122").Value = 53588676
$ws_LTW.Range("L122").Value = 10689.375
$ws_LTW.Range("M122").Value = -53586226
$ws_LTW.Range("N122").Value = -15589.375
$ws_WVR.Range("H62").Value = 3928.9
$ws_WVR.Range("I62").Value = 3464.8333
$ws_WVR.Range("K62").Value = 3464.8333
$ws_WVR.Range("M62").Value = -2840.8333
$ws_WVR.Range("H65").Value = 3928.9
$ws_WVR.Range("I65").Value = 3464.8333
$ws_WVR.Range("K65").Value = 17324.1665
$ws_WVR.Range("M65").Value = -14204.1665
$ws_WVR.Range("H81").Value = 74502.78999999999
$ws_WVR.Range("I81").Value = 169999.83
$ws_WVR.Range("K81").Value = 339999.66
$ws_WVR.Range("M81").Value = -338938.66
$ws_WVR.Range("H84").Value = 74502.78999999999
$ws_WVR.Range("I84").Value = 169999.83
$ws_WVR.Range("K84").Value = 1699998.3
$ws_WVR.Range("M84").Value = -1694694.3
$ws_WVR.Range("H107").Value = 380.69696
$ws_WVR.Range("I107").Value = 263.58334
$ws_WVR.Range("K107").Value = 790.7500200000001
$ws_WVR.Range("M107").Value = 1129.24998
$ws_WVR.Range("H122").Value = 27779574
$ws_WVR.Range("I122").Value = 35715010
$ws_WVR.Range("J122").Value = 5555
$ws_WVR.Range("K122").Value = 107145030
$ws_WVR.Range("L122").Value = 16665
$ws_WVR.Range("M122").Value = -107142580
$ws_WVR.Range("N122").Value = -21565
$ws_WVR.Range("H136").Value = 1290.6279
$ws_WVR.Range("I136").Value = 1344.7222
$ws_WVR.Range("J136").Value = 1012.4286
$ws_WVR.Range("K136").Value = 4034.1666
$ws_WVR.Range("L136").Value = 3037.2858
$ws_WVR.Range("M136").Value = -1484.1666
$ws_WVR.Range("N136").Value = -8137.2858

Write-Output "Applied 193 cell updates across 7 sheets"